$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory labels (column H) to pluralized/renamed forms
$ws.Range("H11").Value = "scatter plot(s)"
$ws.Range("H12").Value = "scatter plot(s)"
$ws.Range("H19").Value = "line graph(s)"
$ws.Range("H20").Value = "line graph(s)"
$ws.Range("H21").Value = "line graph(s)"
$ws.Range("H22").Value = "scatter plot(s)"
$ws.Range("H23").Value = "line graph(s)"
$ws.Range("H26").Value = "line graph(s)"
$ws.Range("H29").Value = "data display"
$ws.Range("H31").Value = "bar chart(s)"
$ws.Range("H36").Value = "line graph(s)"

# Remove the "is_viewed" column (column I) entirely
$ws.Range("I1:I36").EntireColumn.Delete()
